# Fixed a bug in respin
# The reel-weight table (rows 2-21, columns A-F) was being written out in the
# wrong row order. This re-applies the correct order by permuting the rows
# back to how they should be, without altering any of the underlying values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (buggy) row order for the data block A2:F21 before
# writing anything back, so the reorder is computed from a consistent source.
$original = $ws.Range("A2:F21").Value2

# Maps each corrected row (2-21) to the source row it should come from in the
# snapshot above (1-based offset within the 20-row block, i.e. row 2 -> index 1).
$rowMap = @{
    2  = 9
    3  = 2
    4  = 12
    5  = 8
    6  = 4
    7  = 5
    8  = 3
    9  = 10
    10 = 14
    11 = 13
    12 = 7
    13 = 11
    14 = 15
    15 = 6
    16 = 17
    17 = 18
    18 = 16
    19 = 20
    20 = 21
    21 = 19
}

$corrected = New-Object 'object[,]' 20,6

foreach ($destRow in 2..21) {
    $srcRow = $rowMap[$destRow]
    $srcIndex = $srcRow - 1     # 1-based index into $original (row 2 -> 1, row 21 -> 20)
    $destIndex = $destRow - 2   # 0-based index into $corrected (row 2 -> 0, row 21 -> 19)

    for ($col = 1; $col -le 6; $col++) {
        $corrected[$destIndex, $col - 1] = $original[$srcIndex, $col]
    }
}

$ws.Range("A2:F21").Value2 = $corrected
